$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 874.6
$ws.Range("I19").Value = 850.7273
$ws.Range("J19").Value = 893.3570999999999
$ws.Range("K19").Value = 850.7273
$ws.Range("L19").Value = 893.3570999999999
$ws.Range("M19").Value = -675.7273
$ws.Range("N19").Value = -1243.3571

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H38").Value = 1000
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H61").Value = 265335.75
$ws.Range("I61").Value = 1825.8572
$ws.Range("J61").Value = 590848
$ws.Range("K61").Value = 1825.8572
$ws.Range("L61").Value = 590848
$ws.Range("M61").Value = -1613.8572
$ws.Range("N61").Value = -591272
$ws.Range("H101").Value = 42333
$ws.Range("J101").Value = 42333
$ws.Range("L101").Value = 42333
$ws.Range("N101").Value = -48823
$ws.Range("H132").Value = 3867.725
$ws.Range("I132").Value = 850.4483
$ws.Range("J132").Value = 11822.363
$ws.Range("K132").Value = 2551.3449
$ws.Range("L132").Value = 35467.089
$ws.Range("M132").Value = -21.34490000000005
$ws.Range("N132").Value = -40527.089
$ws.Range("H136").Value = 265335.75
$ws.Range("I136").Value = 1825.8572
$ws.Range("J136").Value = 590848
$ws.Range("K136").Value = 5477.571599999999
$ws.Range("L136").Value = 1772544
$ws.Range("M136").Value = -2927.571599999999
$ws.Range("N136").Value = -1777644

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 1008.5455
$ws.Range("I5").Value = 686.125
$ws.Range("K5").Value = 686.125
$ws.Range("M5").Value = -573.125
$ws.Range("H105").Value = 1410.6875
$ws.Range("I105").Value = 1112.8182
$ws.Range("J105").Value = 2066
$ws.Range("K105").Value = 1112.8182
$ws.Range("L105").Value = 2066
$ws.Range("M105").Value = 634.1818000000001
$ws.Range("N105").Value = -5560
$ws.Range("H132").Value = 77900
$ws.Range("J132").Value = 77900
$ws.Range("L132").Value = 77900
$ws.Range("N132").Value = -88020
$ws.Range("H134").Value = 2136.926
$ws.Range("I134").Value = 1267.9412
$ws.Range("K134").Value = 3803.8236
$ws.Range("M134").Value = -1268.8236

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 2927.5
$ws.Range("J13").Value = 2927.5
$ws.Range("L13").Value = 2927.5
$ws.Range("N13").Value = -3205.5
$ws.Range("H19").Value = 338.57144
$ws.Range("I19").Value = 60
$ws.Range("J19").Value = 840
$ws.Range("K19").Value = 60
$ws.Range("L19").Value = 840
$ws.Range("M19").Value = 110
$ws.Range("N19").Value = -1180
$ws.Range("H24").Value = 338.57144
$ws.Range("I24").Value = 60
$ws.Range("J24").Value = 840
$ws.Range("K24").Value = 60
$ws.Range("L24").Value = 840
$ws.Range("M24").Value = 110
$ws.Range("N24").Value = -1180
$ws.Range("H31").Value = 19256062
$ws.Range("I31").Value = 333334080
$ws.Range("J31").Value = 26797.348
$ws.Range("K31").Value = 333334080
$ws.Range("L31").Value = 26797.348
$ws.Range("M31").Value = -333333785
$ws.Range("N31").Value = -27387.348
$ws.Range("H32").Value = 980
$ws.Range("I32").Value = 980
$ws.Range("K32").Value = 980
$ws.Range("M32").Value = -664
$ws.Range("H34").Value = 19256062
$ws.Range("I34").Value = 333334080
$ws.Range("J34").Value = 26797.348
$ws.Range("K34").Value = 333334080
$ws.Range("L34").Value = 26797.348
$ws.Range("M34").Value = -333333878
$ws.Range("N34").Value = -27201.348
$ws.Range("H35").Value = 749.8570999999999
$ws.Range("I35").Value = 749.8570999999999
$ws.Range("K35").Value = 749.8570999999999
$ws.Range("M35").Value = -455.8570999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 768.2857
$ws.Range("I5").Value = 466.77777
$ws.Range("J5").Value = 1785.875
$ws.Range("K5").Value = 1400.33331
$ws.Range("L5").Value = 5357.625
$ws.Range("M5").Value = -1288.33331
$ws.Range("N5").Value = -5581.625
$ws.Range("H44").Value = 55555824
$ws.Range("I44").Value = 201
$ws.Range("J44").Value = 111111450
$ws.Range("K44").Value = 603
$ws.Range("L44").Value = 333334350
$ws.Range("M44").Value = -205
$ws.Range("N44").Value = -333335146
$ws.Range("H135").Value = 768.2857
$ws.Range("I135").Value = 466.77777
$ws.Range("J135").Value = 1785.875
$ws.Range("K135").Value = 4200.99993
$ws.Range("L135").Value = 16072.875
$ws.Range("M135").Value = -1665.99993
$ws.Range("N135").Value = -21142.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3175
$ws.Range("I80").Value = 2200
$ws.Range("K80").Value = 2200
$ws.Range("M80").Value = -1202
$ws.Range("H83").Value = 3175
$ws.Range("I83").Value = 2200
$ws.Range("K83").Value = 11000
$ws.Range("M83").Value = -6008

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 280
$ws.Range("I22").Value = 207.5
$ws.Range("J22").Value = 425
$ws.Range("K22").Value = 207.5
$ws.Range("L22").Value = 425
$ws.Range("M22").Value = 87.5
$ws.Range("N22").Value = -1015
$ws.Range("H27").Value = 280
$ws.Range("I27").Value = 207.5
$ws.Range("J27").Value = 425
$ws.Range("K27").Value = 207.5
$ws.Range("L27").Value = 425
$ws.Range("M27").Value = -100.5
$ws.Range("N27").Value = -639
$ws.Range("H40").Value = 1327.0714
$ws.Range("I40").Value = 1222.8235
$ws.Range("J40").Value = 1488.1818
$ws.Range("K40").Value = 1222.8235
$ws.Range("L40").Value = 1488.1818
$ws.Range("M40").Value = -1086.8235
$ws.Range("N40").Value = -1760.1818

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 3017.1667
$ws.Range("I12").Value = 2967.3333
$ws.Range("J12").Value = 3166.6667
$ws.Range("K12").Value = 2967.3333
$ws.Range("L12").Value = 3166.6667
$ws.Range("M12").Value = -2825.3333
$ws.Range("N12").Value = -3450.6667
$ws.Range("H132").Value = 1897.5636
$ws.Range("I132").Value = 1890.2285
$ws.Range("J132").Value = 1910.4
$ws.Range("K132").Value = 5670.6855
$ws.Range("L132").Value = 5731.200000000001
$ws.Range("M132").Value = -3140.6855
$ws.Range("N132").Value = -10791.2
